$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.97
$summary.Range("B4").Value = -0.03
$summary.Range("B5").Value = -0.04
$summary.Range("B6").Value = 14
$summary.Range("B8").Value = 7
$summary.Range("B9").Value = 35.71

$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.97
$status.Range("D4").Value = 14
$status.Range("E4").Value = -0.03
$status.Range("F4").Value = -0.03
$status.Range("G4").Value = 35.71

$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G15").Value = 0.92
$allTrades.Range("H15").Value = "CLOSED"
$allTrades.Range("I15").Value = -1.0753
$allTrades.Range("J15").Value = -0.01
$allTrades.Range("K15").Value = 99.97
$allTrades.Range("P15").Value = "early_exit"
$allTrades.Range("Q15").Value = 0.13

$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G15").Value = 0.92
$mm.Range("H15").Value = "CLOSED"
$mm.Range("I15").Value = -1.0753
$mm.Range("J15").Value = -0.01
$mm.Range("K15").Value = 99.97
$mm.Range("P15").Value = "early_exit"
$mm.Range("Q15").Value = 0.13
